$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..60 down to 8..60
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44462
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112026
$ws.Cells.Item(7, 7).Value = "Haba"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 85
$ws.Cells.Item(7, 11).Value = 11000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 11529
$ws.Cells.Item(7, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 461
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
